# Repair Variables in view file specifications
#
# The document accumulated a bunch of spots where a single logical run of
# text had been split into two (or more) <w:r> elements that carry
# identical run formatting. This script uses Find/Replace across the
# whole-document text (which is agnostic to the underlying run
# boundaries) so Word rebuilds each spot as a single run, and also fixes
# two real typos ("so-and-so" -> "soAndSo", "criationDate" -> "creationDate").
#
# wdReplaceAll = 2, wdFindContinue = 1

$d = $word.ActiveDocument

function Fix-Text($old, $new) {
    $r = $d.Content
    $ok = $r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        Write-Output "NOT FOUND: $old"
    }
}

# 1) "so-and-s" + "o"  ->  "soAndSo"
Fix-Text "so-and-so" "soAndSo"

# 2) "Evento" + ".php" -> "Evento.php"  (class heading line)
Fix-Text "Evento.php" "Evento.php"

# 3) "Event" + ".php" -> "Event.php"
Fix-Text "Event.php" "Event.php"

# 4) "Creates a “Evento" + ".php” object"
Fix-Text "Creates a “Evento.php” object" "Creates a “Evento.php” object"

# 5) "Receives event's" + " starting date"
Fix-Text "Receives event’s starting date" "Receives event’s starting date"

# 6) "Receives event's" + " ending date"
Fix-Text "Receives event’s ending date" "Receives event’s ending date"

# 7) "Receives event's" + " men price"
Fix-Text "Receives event’s men price" "Receives event’s men price"

# 8) "Receives event's " + "women price"
Fix-Text "Receives event’s women price" "Receives event’s women price"

# 9) "Receives event's " + "promoter"
Fix-Text "Receives event’s promoter" "Receives event’s promoter"

# 10) "Receives event's " + "Facebook event page"
Fix-Text "Receives event’s Facebook event page" "Receives event’s Facebook event page"

# 11) "Receives event's " + "creation date"
Fix-Text "Receives event’s creation date" "Receives event’s creation date"

# 12) "criationDate" (typo) -> "creationDate", split as cr / e / _GoBack / ationDate
$rng = $d.Content
$rng.Find.Execute("criationDate") | Out-Null
$start = $rng.Start

# fix the typo: 3rd character 'i' -> 'e'
$iChar = $d.Range($start + 2, $start + 3)
$iChar.Text = "e"

# force a genuine run boundary between "cr" and "e" (identical formatting on
# both sides, so a no-op bold toggle is used purely to split the run)
$crRun = $d.Range($start, $start + 2)
$crRun.Font.Bold = 1
$crRun.Font.Bold = 0

# move (or create) the _GoBack bookmark between "e" and "ationDate"; Word
# keeps a single _GoBack so this also removes it from its previous location
$bmRange = $d.Range($start + 3, $start + 3)
$d.Bookmarks.Add("_GoBack", $bmRange)

# 13) "Receives event's " + "description"
Fix-Text "Receives event’s description" "Receives event’s description"

# 14) "Receives event's " + "tickets lot"
Fix-Text "Receives event’s tickets lot" "Receives event’s tickets lot"

# 15) "Receives event's " + "thumbnail"
Fix-Text "Receives event’s thumbnail" "Receives event’s thumbnail"

# 16) "Receives event's " + "age recommendation"
Fix-Text "Receives event’s age recommendation" "Receives event’s age recommendation"

# 17) "editarPessoa" + ".php"
Fix-Text "editarPessoa.php" "editarPessoa.php"

# 18) "Contains Person first name" + " value"
Fix-Text "Contains Person first name value" "Contains Person first name value"

# 19) "Contains Person email" + " value"
Fix-Text "Contains Person email value" "Contains Person email value"

# 20) "Contains Person " + "ID" + " value"
Fix-Text "Contains Person ID value" "Contains Person ID value"

# 21) "Contains Person " + "phone number"
Fix-Text "Contains Person phone number" "Contains Person phone number"

Write-Output "done"
